$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H8").Value = "May 24th, 2018"
$ws.Range("H3").Select()
